$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: registered families (E4:K4)
$ws.Range("E4").Value = 8322
$ws.Range("F4").Value = 5355
$ws.Range("G4").Value = 5161
$ws.Range("H4").Value = 4888
$ws.Range("I4").Value = 4781
$ws.Range("J4").Value = 5136
$ws.Range("K4").Value = 5425

# Row 5: subsistence allowance recipient families (E5:K5)
$ws.Range("E5").Value = 2562
$ws.Range("F5").Value = 2485
$ws.Range("G5").Value = 2231
$ws.Range("H5").Value = 1987
$ws.Range("I5").Value = 1750
$ws.Range("J5").Value = 2343
$ws.Range("K5").Value = 2710
